# Update cryptocurrency price/volume/time data to reflect the latest
# snapshot fetched by the GitHub Actions symbol-list updater.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    # Force the cell to stay text (these numeric-looking strings must
    # not be reinterpreted as numbers/percentages by Excel).
    $range.NumberFormat = "@"
    $range.Value = $text
    # Restore the default style so formatting matches the original file.
    $range.Style = "Normal"
}

Set-TextValue "D2" "289.20"
Set-TextValue "E2" "-0.31%"
Set-TextValue "G2" "15"

Set-TextValue "D3" "31.04"
Set-TextValue "E3" "2.33%"
Set-TextValue "G3" "15"

Set-TextValue "D4" "4.949"
Set-TextValue "E4" "0.70%"
Set-TextValue "G4" "15"

Set-TextValue "D5" "0.07369"
Set-TextValue "E5" "1.74%"
Set-TextValue "G5" "15"

Set-TextValue "D6" "2.290"
Set-TextValue "E6" "27.35%"
Set-TextValue "G6" "15"

Set-TextValue "D7" "7.682"
Set-TextValue "E7" "1.13%"
Set-TextValue "G7" "15"

Set-TextValue "D8" "0.9174"
Set-TextValue "E8" "1.83%"
Set-TextValue "G8" "15"

Set-TextValue "D9" "0.09238"
Set-TextValue "E9" "16.58%"
Set-TextValue "G9" "15"

Set-TextValue "D10" "0.1702"
Set-TextValue "E10" "1.81%"
Set-TextValue "G10" "15"

Set-TextValue "D11" "0.08269"
Set-TextValue "E11" "2.78%"
Set-TextValue "G11" "15"

Set-TextValue "D12" "0.03104"
Set-TextValue "E12" "2.30%"
Set-TextValue "G12" "15"

Set-TextValue "D13" "0.09988"
Set-TextValue "E13" "-0.37%"
Set-TextValue "G13" "15"

Set-TextValue "D14" "0.001495"
Set-TextValue "E14" "0.00%"
Set-TextValue "G14" "15"

Set-TextValue "D15" "0.005719"
Set-TextValue "E15" "-0.08%"
Set-TextValue "G15" "15"

Set-TextValue "D16" "3.469"
Set-TextValue "E16" "-0.35%"
Set-TextValue "G16" "15"

Set-TextValue "D17" "3.735"
Set-TextValue "E17" "0.83%"
Set-TextValue "G17" "15"

Set-TextValue "D18" "2.124"
Set-TextValue "E18" "2.30%"
Set-TextValue "G18" "15"

Set-TextValue "D19" "0.3325"
Set-TextValue "E19" "0.21%"
Set-TextValue "G19" "15"

Set-TextValue "D20" "0.1290"
Set-TextValue "E20" "-0.97%"
Set-TextValue "G20" "15"

Set-TextValue "D21" "4.148"
Set-TextValue "E21" "4.81%"
Set-TextValue "G21" "15"

Set-TextValue "D22" "0.2122"
Set-TextValue "E22" "-2.12%"
Set-TextValue "G22" "15"

Set-TextValue "D23" "0.04506"
Set-TextValue "E23" "0.03%"
Set-TextValue "G23" "15"

Set-TextValue "D24" "0.001215"
Set-TextValue "E24" "0.22%"
Set-TextValue "G24" "15"

Set-TextValue "D25" "0.004196"
Set-TextValue "E25" "-5.34%"
Set-TextValue "G25" "15"

Set-TextValue "D26" "0.0001298"
Set-TextValue "E26" "0.14%"
Set-TextValue "G26" "15"

Set-TextValue "D27" "0.0003389"
Set-TextValue "G27" "15"

Set-TextValue "G28" "15"

Set-TextValue "G29" "15"

Set-TextValue "G30" "15"

Set-TextValue "G31" "15"

Set-TextValue "G32" "15"

Set-TextValue "G33" "15"

Set-TextValue "G34" "15"

Set-TextValue "G35" "15"

Set-TextValue "G36" "15"

Set-TextValue "G37" "15"

Set-TextValue "G38" "15"

Set-TextValue "D39" "0.01583"
Set-TextValue "E39" "0.57%"
Set-TextValue "G39" "15"

Set-TextValue "D40" "0.04494"
Set-TextValue "E40" "3.29%"
Set-TextValue "G40" "15"

Set-TextValue "D41" "0.007375"
Set-TextValue "E41" "0.90%"
Set-TextValue "G41" "15"

Set-TextValue "D42" "0.009842"
Set-TextValue "E42" "-2.08%"
Set-TextValue "G42" "15"

Set-TextValue "D43" "0.1336"
Set-TextValue "E43" "1.91%"
Set-TextValue "G43" "15"

Set-TextValue "D44" "0.002227"
Set-TextValue "E44" "8.93%"
Set-TextValue "G44" "15"

Set-TextValue "D45" "0.008938"
Set-TextValue "E45" "-1.34%"
Set-TextValue "G45" "15"

Set-TextValue "D46" "0.00006095"
Set-TextValue "E46" "3.53%"
Set-TextValue "G46" "15"

Set-TextValue "D47" "0.00000000749"
Set-TextValue "E47" "-0.06%"
Set-TextValue "G47" "15"

Set-TextValue "G48" "15"

Set-TextValue "D49" "0.002097"
Set-TextValue "G49" "15"

Set-TextValue "D50" "0.00002097"
Set-TextValue "E50" "-0.06%"
Set-TextValue "G50" "15"

Set-TextValue "D51" "0.0001997"
Set-TextValue "E51" "-0.06%"
Set-TextValue "G51" "15"
